$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add new column U ---
$ws.Range("U1").Value = "Fertilizer Recommendation"

# --- Row 2 updates ---

# B2 holds a date-like literal that must stay as plain text "11-04-2024"
# (not get auto-converted into a date serial number). Temporarily force
# Text format, assign the value, then drop back to the default "Normal"
# style so no stray number formatting is left behind on the cell.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "11-04-2024"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").Value = 546
$ws.Range("D2").Value = 545
$ws.Range("E2").Value = "asdjhalksdjkl"
$ws.Range("F2").Value = 36
$ws.Range("G2").Value = "Female"
$ws.Range("H2").Value = 54
$ws.Range("I2").Value = "214asd45asd65"

# J2 is a purely-numeric mobile number that must remain text, otherwise it
# gets auto-converted to a numeric value. Same Text-format trick as B2.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "8797865461"
$ws.Range("J2").Style = "Normal"

$ws.Range("K2").Value = 7.5
$ws.Range("L2").Value = 154
$ws.Range("M2").Value = 123
$ws.Range("N2").Value = 212
$ws.Range("P2").Value = 21
$ws.Range("Q2").Value = 56
$ws.Range("R2").Value = 65
$ws.Range("S2").Value = 0.568192230958883
$ws.Range("T2").Value = "Grow maize, soybean, groundnut, cotton, and incorporate legumes into the cropping system."
$ws.Range("U2").Value = "Apply organic amendments like Compost (2-3 tonnes/ha), Vermicompost (1-1.5 tonnes/ha), or Well-decomposed Farmyard manure (5-7.5 tonnes/ha). Use biofertilizers like Rhizobium (200-300 g/ha), Azotobacter (200-300 g/ha), and PSB (500-750 g/ha). Apply chemical fertilizers at 75% of the recommended dose based on soil test results and crop requirements."
